$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the table (order matters for shared string indices)
$ws.Range("I6").Value = "Total"
$ws.Range("J6").Value = "Column2"
$ws.Range("F6").Value = "Parameter"

# Column J width (pixel width ~72px @ Calibri 11 / MDW7 -> closest achievable stored width)
$ws.Columns.Item(10).ColumnWidth = 9.451822916666666

# Create the Excel Table (ListObject) over F6:J29
$tbl = $ws.ListObjects.Add(1, $ws.Range("F6:J29"), [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# View changes: zoom + selection
$aw = $excel.ActiveWindow
$aw.Zoom = 130
[void]$ws.Range("D14").Select()

Write-Host "done"
